# Insert a new column "ID" at the start of the sheet (A), shifting the
# existing columns A:E (A,B,C,D,F headers) to B:F. Then populate the new
# column A with a header and per-row sample identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data one column to the right.
$ws.Columns("A").Insert()

# Give the new header cell the same look as the rest of the header row
# (bold, bordered, centered), then set its text.
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "ID"

# New row labels for column A.
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
